$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.725.36"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.08%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.290.78"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.98%  "
$ws.Range("E4").Value = "  -0.24%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "102.69"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +6.30%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "270.10"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.06%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.619"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.06%  "
$ws.Range("E8").Value = "  -0.16%  "
$ws.Range("E9").Value = "  -2.06%  "
$ws.Range("E10").Value = "  -0.19%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0935"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.01%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.98"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.34%  "
$ws.Range("E13").Value = "  +1.69%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "15.77"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.20%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.857"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.09%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.305.03"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.44%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "43.713.06"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.03%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0000110"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.27%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.26"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.80%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "72.27"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.46%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.50"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +10.09%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "233.31"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.36%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.87"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +14.10%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.16"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.15%  "
$ws.Range("E25").Value = "  +0.03%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.22"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.76%  "
$ws.Range("E27").Value = "  -0.26%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "39.41"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.97%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.23"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.01%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "177.33"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.62%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "21.78"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.63%  "
$ws.Range("E32").Value = "  +0.28%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.46"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.10%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.89"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +12.27%  "
$ws.Range("E35").Value = "  +0.13%  "
$ws.Range("E36").Value = "  +0.85%  "
$ws.Range("E37").Value = "  -2.06%  "
$ws.Range("E38").Value = "  +7.05%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.33"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.39%  "
$ws.Range("E40").Value = "  -3.56%  "
$ws.Range("E41").Value = "  +1.83%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "12.25"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.31%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "65.36"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +5.79%  "
$ws.Range("E44").Value = "  -1.67%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.77"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.98%  "
$ws.Range("E46").Value = "  -0.90%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.22"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.12%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "98.48"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.62%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.452"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +9.37%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.54"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +12.40%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.513.87"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.18%  "
